$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW35-FE-LIFTER ---
$ws1 = $wb.Worksheets.Item(1)
$r = 55
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item($r - 1, 1).NumberFormat
$ws1.Cells.Item($r, 1).Value = 45751.85768704861
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x72"
$ws1.Cells.Item($r, 5).Value = "0xd"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r, 8).Value = 370
$ws1.Cells.Item($r, 9).Value = 13

# --- Sheet 2: ROW35-MID-LIFTER ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item($r - 1, 1).NumberFormat
$ws2.Cells.Item($r, 1).Value = 45751.71048190972
$ws2.Cells.Item($r, 2).Value = "0x01,0x90"
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws2.Cells.Item($r, 5).Value = "0xe"
$ws2.Cells.Item($r, 6).Value = 400
$ws2.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws2.Cells.Item($r, 8).Value = 366
$ws2.Cells.Item($r, 9).Value = 14

# --- Sheet 3: ROW02-FE-LIFTER ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item($r, 1).NumberFormat = $ws3.Cells.Item($r - 1, 1).NumberFormat
$ws3.Cells.Item($r, 1).Value = 45751.85300739583
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x72"
$ws3.Cells.Item($r, 5).Value = "0x3"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r, 8).Value = 370
$ws3.Cells.Item($r, 9).Value = 3

# --- Sheet 4: ROW02-MID-LIFTER ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item($r, 1).NumberFormat = $ws4.Cells.Item($r - 1, 1).NumberFormat
$ws4.Cells.Item($r, 1).Value = 45751.9083777662
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x6e"
$ws4.Cells.Item($r, 5).Value = "0x3"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = [double]"9.85046333984776e+23"
$ws4.Cells.Item($r, 8).Value = 366
$ws4.Cells.Item($r, 9).Value = 3
